$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numRange = $ws.Range('D5,D6,D9,D10,D12,D16,D18,D19,D20,D21,D22,D23,D26,D28,D29,D31,D32,D34,D36,D37,D38,D39,D41,D43,D44,D45,D48,D49,D50,D51')
$numRange.NumberFormat = "@"

$ws.Range("D2").Value = '60.865.80'
$ws.Range("E2").Value = '  -5.01%  '
$ws.Range("D3").Value = '3.272.38'
$ws.Range("E3").Value = '  -5.77%  '
$ws.Range("D5").Value = '562.36'
$ws.Range("E5").Value = '  -3.86%  '
$ws.Range("D6").Value = '126.29'
$ws.Range("E6").Value = '  -4.02%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.273.49'
$ws.Range("E8").Value = '  -5.71%  '
$ws.Range("D9").Value = '0.473'
$ws.Range("E9").Value = '  -1.92%  '
$ws.Range("D10").Value = '7.25'
$ws.Range("E10").Value = '  -5.08%  '
$ws.Range("E11").Value = '  -4.71%  '
$ws.Range("D12").Value = '0.371'
$ws.Range("E12").Value = '  -3.99%  '
$ws.Range("D13").Value = '3.839.77'
$ws.Range("E13").Value = '  -5.58%  '
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("D15").Value = '3.290.29'
$ws.Range("E15").Value = '  -5.32%  '
$ws.Range("D16").Value = '0.0000165'
$ws.Range("E16").Value = '  -6.57%  '
$ws.Range("D17").Value = '60.991.24'
$ws.Range("E17").Value = '  -4.82%  '
$ws.Range("D18").Value = '23.85'
$ws.Range("E18").Value = '  -1.97%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '13.14'
$ws.Range("E19").Value = '  -2.32%  '
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '5.54'
$ws.Range("E20").Value = '  -2.91%  '
$ws.Range("D21").Value = '8.82'
$ws.Range("E21").Value = '  -11.44%  '
$ws.Range("D22").Value = '349.41'
$ws.Range("E22").Value = '  -9.08%  '
$ws.Range("D23").Value = '0.547'
$ws.Range("E23").Value = '  -4.73%  '
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").Value = '3.412.11'
$ws.Range("E25").Value = '  -5.55%  '
$ws.Range("D26").Value = '68.77'
$ws.Range("E26").Value = '  -7.90%  '
$ws.Range("E27").Value = '  -5.75%  '
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("D29").Value = '7.00'
$ws.Range("E29").Value = '  -2.18%  '
$ws.Range("E30").Value = '  -0.78%  '
$ws.Range("D31").Value = '2.09'
$ws.Range("E31").Value = '  -6.38%  '
$ws.Range("D32").Value = '7.70'
$ws.Range("E32").Value = '  -2.89%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").Value = '0.147'
$ws.Range("E34").Value = '  -3.29%  '
$ws.Range("D35").Value = '3.306.12'
$ws.Range("E35").Value = '  -5.61%  '
$ws.Range("D36").Value = '22.40'
$ws.Range("E36").Value = '  -2.32%  '
$ws.Range("D37").Value = '5.18'
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").Value = '6.71'
$ws.Range("E38").Value = '  -0.82%  '
$ws.Range("D39").Value = '162.96'
$ws.Range("E39").Value = '  +0.44%  '
$ws.Range("E40").Value = '  -3.68%  '
$ws.Range("D41").Value = '0.0745'
$ws.Range("E41").Value = '  -3.88%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("D43").Value = '40.88'
$ws.Range("E43").Value = '  -1.29%  '
$ws.Range("D44").Value = '4.30'
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("D45").Value = '0.733'
$ws.Range("E45").Value = '  -8.21%  '
$ws.Range("E46").Value = '  -1.54%  '
$ws.Range("E47").Value = '  -5.38%  '
$ws.Range("D48").Value = '21.95'
$ws.Range("E48").Value = '  -7.80%  '
$ws.Range("D49").Value = '6.62'
$ws.Range("E49").Value = '  -1.32%  '
$ws.Range("D50").Value = '0.843'
$ws.Range("E50").Value = '  -7.62%  '
$ws.Range("D51").Value = '20.83'
$ws.Range("E51").Value = '  +1.95%  '

$numRange.ClearFormats()
